$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 Status/Date changed: order ORD-1736840768506-3f916962 now shows
# "return approved" and keeps its original timestamp.
$ws.Cells.Item(19, 8).Value = "return approved"
$ws.Cells.Item(19, 9).Value = "14/01/2025, 13:16:08"

# New sales-report rows (20-32) appended below the existing data.
$newRows = @(
    @{ Row = 20; A = "ORD-1736921273602-1d81a9e1"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 783; E = 196; F = "NEWYEAR25";  G = 587;  H = "pending";         I = "15/01/2025, 11:37:53" },
    @{ Row = 21; A = "ORD-1736922318507-df33f9bd"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 178; F = "NEWYEAR25";  G = 534;  H = "pending";         I = "15/01/2025, 11:55:18" },
    @{ Row = 22; A = "ORD-1736922879834-aa9a9186"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 128; F = "MAX50";      G = 712;  H = "pending";         I = "15/01/2025, 12:04:39" },
    @{ Row = 23; A = "ORD-1736924311922-71beabcc"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "pending";         I = "15/01/2025, 12:28:31" },
    @{ Row = 24; A = "ORD-1736924400415-f87fdcf3"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 783; E = 141; F = "MAX50";      G = 642;  H = "canceled";        I = "15/01/2025, 12:30:00" },
    @{ Row = 25; A = "ORD-1736950818225-ddb68a48"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "canceled";        I = "15/01/2025, 19:50:18" },
    @{ Row = 26; A = "ORD-1736950882256-eda8eb06"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "return approved"; I = "15/01/2025, 19:51:22" },
    @{ Row = 27; A = "ORD-1737001643934-dcd640af"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 391; E = 0;   F = "null";       G = 391;  H = "pending";         I = "16/01/2025, 09:57:23" },
    @{ Row = 28; A = "ORD-1737001988963-2000ff84"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 783; E = 0;   F = "null";       G = 783;  H = "pending";         I = "16/01/2025, 10:03:08" },
    @{ Row = 29; A = "ORD-1737005141372-6f847034"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "pending";         I = "16/01/2025, 10:55:41" },
    @{ Row = 30; A = "ORD-1737005481173-86d4f49c"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "pending";         I = "16/01/2025, 11:01:21" },
    @{ Row = 31; A = "ORD-1737006109565-f7dead73"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 712; E = 0;   F = "null";       G = 712;  H = "canceled";        I = "16/01/2025, 11:11:49" },
    @{ Row = 32; A = "ORD-1737006560056-2fbc8a89"; B = "Prithviraj"; C = "iamprithvi@gmail.com"; D = 783; E = 0;   F = "null";       G = 783;  H = "pending";         I = "16/01/2025, 11:19:20" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}
